$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "96.617.59"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "3.664.87"
$ws.Range("E3").Value = "  +1.74%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue "D5" "239.51"
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("E6").Value = "  +10.15%  "

Set-TextValue "D7" "655.47"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("E8").Value = "  +0.28%  "

$ws.Range("E9").Value = "  +2.50%  "

Set-TextValue "D10" "1.00"
$ws.Range("E10").Value = "  +0.11%  "

$ws.Range("D11").Value = "3.661.31"
$ws.Range("E11").Value = "  +1.62%  "

Set-TextValue "D12" "45.22"
$ws.Range("E12").Value = "  +2.48%  "

Set-TextValue "D13" "0.205"
$ws.Range("E13").Value = "  +0.18%  "

Set-TextValue "D14" "6.78"
$ws.Range("E14").Value = "  +4.74%  "

$ws.Range("D15").Value = "4.348.62"
$ws.Range("E15").Value = "  +1.80%  "

Set-TextValue "D16" "0.0000271"
$ws.Range("E16").Value = "  +3.95%  "

$ws.Range("D17").Value = "96.389.60"
$ws.Range("E17").Value = "  -1.38%  "

Set-TextValue "D18" "8.82"
$ws.Range("E18").Value = "  +13.57%  "

$ws.Range("D19").Value = "3.664.50"
$ws.Range("E19").Value = "  +1.81%  "

Set-TextValue "D20" "18.73"
$ws.Range("E20").Value = "  +4.16%  "

Set-TextValue "D21" "12.70"
$ws.Range("E21").Value = "  -0.02%  "

Set-TextValue "D22" "0.525"
$ws.Range("E22").Value = "  +0.86%  "

Set-TextValue "D23" "527.37"
$ws.Range("E23").Value = "  +1.87%  "

Set-TextValue "D24" "3.49"
$ws.Range("E24").Value = "  -1.06%  "

Set-TextValue "D25" "7.08"
$ws.Range("E25").Value = "  +2.86%  "

Set-TextValue "D26" "0.0000204"
$ws.Range("E26").Value = "  -0.93%  "

Set-TextValue "D27" "102.16"
$ws.Range("E27").Value = "  +0.32%  "

Set-TextValue "D28" "13.52"
$ws.Range("E28").Value = "  +3.92%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D29" "0.167"
$ws.Range("E29").Value = "  +5.98%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D30" "12.41"
$ws.Range("E30").Value = "  +5.24%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "3.03"
$ws.Range("E31").Value = "  +0.39%  "

$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D32" "0.999"
$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D33" "1.92"
$ws.Range("E33").Value = "  +16.71%  "

$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D34" "0.185"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D35" "669.21"
$ws.Range("E35").Value = "  +8.42%  "

$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D36" "1.00"
$ws.Range("E36").Value = "  +0.57%  "

$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D37" "32.41"
$ws.Range("E37").Value = "  +1.86%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D38" "0.595"
$ws.Range("E38").Value = "  +4.25%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D39" "8.83"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D40" "0.160"
$ws.Range("E40").Value = "  +3.95%  "

$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D41" "1.99"
$ws.Range("E41").Value = "  +0.05%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "39.39"
$ws.Range("E42").Value = "  +19.67%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "6.52"
$ws.Range("E43").Value = "  +8.86%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D44" "0.956"
$ws.Range("E44").Value = "  +3.34%  "

$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D45" "1.00"
$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0459"
$ws.Range("E46").Value = "  +4.14%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D47" "0.440"
$ws.Range("E47").Value = "  +14.29%  "

$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
Set-TextValue "D48" "3.83"
$ws.Range("E48").Value = "  +7.14%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D49" "2.31"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D50" "23.66"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "8.64"
$ws.Range("E51").Value = "  +1.02%  "
